$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix JLCPCB part numbers (column D) on the affected BOM rows
$ws.Range("D2").Value = "C15008"
$ws.Range("D3").Value = "C71190"
$ws.Range("D6").Value = "C21189"

# Leave the selection on the last-edited part-number cell
[void]$ws.Range("D13").Select()
